# Demo-version update: append one new record to the "Cases" table and one
# new record to the "Contacts" table (each table grows by exactly one row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Cases table: append row Id=5 (becomes worksheet row 6, A1:J5 -> A1:J6)
# ---------------------------------------------------------------------
$cases    = $wb.Worksheets.Item("Cases")
$casesTbl = $cases.ListObjects.Item("Cases")
$newCase  = $casesTbl.ListRows.Add()
$cr       = $newCase.Range.Row

# Pick up the existing date/datetime number formats (columns B and C) from
# an existing data row so the new cells keep the same display format.
$cases.Range("B2:C2").Copy($cases.Range("B" + $cr + ":C" + $cr))

$cases.Cells.Item($cr, 1).Value = 5                        # Id
$cases.Cells.Item($cr, 2).Value = 44320                    # Test Date
$cases.Cells.Item($cr, 3).Value = 44320.987027581          # Added Date
$cases.Cells.Item($cr, 4).Value = "OX2"                    # Postcode
$cases.Cells.Item($cr, 5).Value = $true                    # Traced?
$cases.Cells.Item($cr, 6).Value = 0                         # Dropped times
$cases.Cells.Item($cr, 7).Value = $false                   # Dropped?
$cases.Cells.Item($cr, 8).Value = "04/05/2021 23:42:46"    # Traced Date

# ---------------------------------------------------------------------
# Contacts table: append row Id=3 (becomes worksheet row 4, A1:F3 -> A1:F4)
# ---------------------------------------------------------------------
$contacts    = $wb.Worksheets.Item("Contacts")
$contactsTbl = $contacts.ListObjects.Item("Contacts")
$newContact  = $contactsTbl.ListRows.Add()
$kr          = $newContact.Range.Row

# Keep the same "Added Date" number format as the rest of the column.
$contacts.Range("C2").Copy($contacts.Range("C" + $kr))

$contacts.Cells.Item($kr, 1).Value = 3                          # Id
$contacts.Cells.Item($kr, 2).Value = 5                          # CaseId (-> new Case row)
$contacts.Cells.Item($kr, 3).Value = 44320.9875444792           # Added Date
$contacts.Cells.Item($kr, 4).Value = "04/05/2021 23:42:46"      # Traced Date
$contacts.Cells.Item($kr, 5).Value = "04/05/2021 23:42:13"      # Contacted date
